$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Value)
    $Cell.Value = "'" + $Value
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '29.015.65'
$ws.Range('E2').Value = '  -0.15%  '
Set-TextValue $ws.Range('D3') '1.832.80'
$ws.Range('E3').Value = '  +0.20%  '
Set-TextValue $ws.Range('D4') '0.9980'
$ws.Range('E4').Value = '  -0.05%  '
Set-TextValue $ws.Range('D5') '244.20'
$ws.Range('E5').Value = '  +1.38%  '
Set-TextValue $ws.Range('D6') '0.6337'
$ws.Range('E6').Value = '  +1.79%  '
Set-TextValue $ws.Range('D7') '0.9998'
$ws.Range('E7').Value = '  -0.03%  '
Set-TextValue $ws.Range('D8') '0.07579'
$ws.Range('E8').Value = '  +2.85%  '
Set-TextValue $ws.Range('D9') '0.2950'
$ws.Range('E9').Value = '  +0.87%  '
Set-TextValue $ws.Range('D10') '22.79'
Set-TextValue $ws.Range('D11') '0.07736'
$ws.Range('E11').Value = '  +0.98%  '
Set-TextValue $ws.Range('D12') '1.836.82'
$ws.Range('E12').Value = '  +0.35%  '
Set-TextValue $ws.Range('D13') '4.994'
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('E14').Value = '  +1.15%  '
Set-TextValue $ws.Range('D15') '83.16'
$ws.Range('E15').Value = '  +1.21%  '
Set-TextValue $ws.Range('D16') '0.000009854'
$ws.Range('E16').Value = '  +9.13%  '
Set-TextValue $ws.Range('D17') '6.115'
$ws.Range('E17').Value = '  +1.36%  '
Set-TextValue $ws.Range('D18') '29.054.72'
$ws.Range('E18').Value = '  +0.01%  '
Set-TextValue $ws.Range('D20') '226.65'
$ws.Range('E20').Value = '  +0.43%  '
Set-TextValue $ws.Range('D21') '0.9990'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('E22').Value = '  +0.84%  '
Set-TextValue $ws.Range('D23') '0.9994'
$ws.Range('E23').Value = '  -0.06%  '
Set-TextValue $ws.Range('D24') '160.33'
$ws.Range('E24').Value = '  +0.47%  '
Set-TextValue $ws.Range('D25') '0.1406'
$ws.Range('E25').Value = '  +3.56%  '
Set-TextValue $ws.Range('D26') '8.534'
$ws.Range('E26').Value = '  +1.34%  '
Set-TextValue $ws.Range('D27') '17.94'
$ws.Range('E27').Value = '  +0.67%  '
Set-TextValue $ws.Range('D28') '1.504'
$ws.Range('E28').Value = '  +0.57%  '
Set-TextValue $ws.Range('D29') '4.121'
$ws.Range('E29').Value = '  +1.61%  '
Set-TextValue $ws.Range('D30') '4.051'
$ws.Range('E30').Value = '  +0.38%  '
Set-TextValue $ws.Range('D31') '1.204'
$ws.Range('E31').Value = '  +0.41%  '
Set-TextValue $ws.Range('D32') '0.05399'
$ws.Range('E32').Value = '  +2.81%  '
Set-TextValue $ws.Range('D33') '1.862'
$ws.Range('E33').Value = '  +0.91%  '
Set-TextValue $ws.Range('D34') '0.7466'
$ws.Range('E34').Value = '  +1.81%  '
Set-TextValue $ws.Range('D35') '1.141'
$ws.Range('E35').Value = '  -1.03%  '
Set-TextValue $ws.Range('D36') '2.667'
$ws.Range('E36').Value = '  +0.74%  '
Set-TextValue $ws.Range('D37') '1.241.69'
$ws.Range('E37').Value = '  -3.95%  '
Set-TextValue $ws.Range('D38') '0.01795'
$ws.Range('E38').Value = '  +0.55%  '
Set-TextValue $ws.Range('D39') '2.757'
$ws.Range('E39').Value = '  +0.34%  '
$ws.Range('E40').Value = '  +4.94%  '
Set-TextValue $ws.Range('D41') '0.9016'
$ws.Range('E41').Value = '  +0.05%  '
Set-TextValue $ws.Range('D42') '1.000'
$ws.Range('E42').Value = '  +0.08%  '
Set-TextValue $ws.Range('D43') '102.62'
$ws.Range('E43').Value = '  +0.80%  '
$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D44') '0.00000000126'
$ws.Range('E44').Value = '  +5.56%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range('D45') '1.985.63'
$ws.Range('E45').Value = '  +0.50%  '
Set-TextValue $ws.Range('D46') '64.84'
$ws.Range('E46').Value = '  +1.26%  '
Set-TextValue $ws.Range('D47') '0.5107'
$ws.Range('E47').Value = '  -0.13%  '
Set-TextValue $ws.Range('D48') '0.4103'
$ws.Range('E48').Value = '  +3.40%  '
Set-TextValue $ws.Range('D49') '9.011'
$ws.Range('E49').Value = '  +2.12%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D50') '6.774'
$ws.Range('E50').Value = '  +1.59%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D51') '0.05779'
$ws.Range('E51').Value = '  +0.07%  '
